$d = $word.ActiveDocument

function Merge-Text($old, $new) {
    $r = $d.Content
    $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Merge-Text "nome nuovo allergene nullo" "nome nuovo allergene nullo"
Merge-Text "nome nuovo allergene uguale ad un allergene già registrato nel sistema" "nome nuovo allergene uguale ad un allergene già registrato nel sistema"
Merge-Text "categoria non specificata" "categoria non specificata"
Merge-Text "categoria non esistente nel sistema" "categoria non esistente nel sistema"
Merge-Text "nome allergene cercato vuoto o nullo" "nome allergene cercato vuoto o nullo"
Merge-Text "nome di un allergene esistente" "nome di un allergene esistente"
Merge-Text "nome di un allergene non esistente" "nome di un allergene non esistente"
Merge-Text "nome dell’allergene da modificare nullo o vuoto" "nome dell’allergene da modificare nullo o vuoto"
Merge-Text "nome allergene vuoto o nullo" "nome allergene vuoto o nullo"

# --- Add new "DAO Categoria" section content ---
$p = $d.Paragraphs.Last
$p.Range.Text = "Elenco categorie presenti nel sistema"

$p2 = $d.Paragraphs.Last
$insertPos = $p2.Range.End - 1
$rng = $d.Range($insertPos, $insertPos)

$xmlFrag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:tbl>
  <w:tblPr>
    <w:tblStyle w:val="TableGrid"/>
    <w:tblW w:w="0" w:type="auto"/>
    <w:tblBorders>
      <w:top w:val="none" w:sz="0" w:space="0" w:color="auto"/>
      <w:left w:val="none" w:sz="0" w:space="0" w:color="auto"/>
      <w:bottom w:val="none" w:sz="0" w:space="0" w:color="auto"/>
      <w:right w:val="none" w:sz="0" w:space="0" w:color="auto"/>
      <w:insideH w:val="none" w:sz="0" w:space="0" w:color="auto"/>
      <w:insideV w:val="none" w:sz="0" w:space="0" w:color="auto"/>
    </w:tblBorders>
    <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
  </w:tblPr>
  <w:tblGrid>
    <w:gridCol w:w="4814"/>
    <w:gridCol w:w="4814"/>
  </w:tblGrid>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4814" w:type="dxa"/>
        <w:tcBorders>
          <w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/>
        </w:tcBorders>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:b/>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:b/>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
          <w:t>dati in ingresso</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4814" w:type="dxa"/>
        <w:tcBorders>
          <w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/>
        </w:tcBorders>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:b/>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:b/>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
          <w:t>risultato atteso</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4814" w:type="dxa"/>
        <w:tcBorders>
          <w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/>
          <w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/>
        </w:tcBorders>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
          <w:lastRenderedPageBreak/>
          <w:t>nessuno</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4814" w:type="dxa"/>
        <w:tcBorders>
          <w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/>
          <w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/>
        </w:tcBorders>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
          <w:t>lista non nulla, con 0 o più elementi</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
</w:tbl>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading2"/>
  </w:pPr>
  <w:r>
    <w:t>Selezione categoria</w:t>
  </w:r>
</w:p>
<w:tbl>
  <w:tblPr>
    <w:tblStyle w:val="TableGrid"/>
    <w:tblW w:w="0" w:type="auto"/>
    <w:tblBorders>
      <w:top w:val="none" w:sz="0" w:space="0" w:color="auto"/>
      <w:left w:val="none" w:sz="0" w:space="0" w:color="auto"/>
      <w:bottom w:val="none" w:sz="0" w:space="0" w:color="auto"/>
      <w:right w:val="none" w:sz="0" w:space="0" w:color="auto"/>
      <w:insideH w:val="none" w:sz="0" w:space="0" w:color="auto"/>
      <w:insideV w:val="none" w:sz="0" w:space="0" w:color="auto"/>
    </w:tblBorders>
    <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
  </w:tblPr>
  <w:tblGrid>
    <w:gridCol w:w="4814"/>
    <w:gridCol w:w="4814"/>
  </w:tblGrid>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4814" w:type="dxa"/>
        <w:tcBorders>
          <w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/>
        </w:tcBorders>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:b/>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:b/>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
          <w:t>dati in ingresso</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4814" w:type="dxa"/>
        <w:tcBorders>
          <w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/>
        </w:tcBorders>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:b/>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:b/>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
          <w:t>risultato atteso</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4814" w:type="dxa"/>
        <w:tcBorders>
          <w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/>
          <w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/>
        </w:tcBorders>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
          <w:t>nome categoria vuoto o nullo</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4814" w:type="dxa"/>
        <w:tcBorders>
          <w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/>
          <w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/>
        </w:tcBorders>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
          <w:t>?</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4814" w:type="dxa"/>
        <w:tcBorders>
          <w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/>
          <w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/>
        </w:tcBorders>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
          <w:t>nome di una categoria registrata nel sistema</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4814" w:type="dxa"/>
        <w:tcBorders>
          <w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/>
          <w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/>
        </w:tcBorders>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
          <w:t>categoria cercata, con elenco allergeni valorizzato</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4814" w:type="dxa"/>
        <w:tcBorders>
          <w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/>
          <w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/>
        </w:tcBorders>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
          <w:t xml:space="preserve">nome di una categoria </w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
          <w:t xml:space="preserve">non </w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
          <w:t>registrata nel sistema</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4814" w:type="dxa"/>
        <w:tcBorders>
          <w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/>
          <w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/>
        </w:tcBorders>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
          <w:t>null</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
</w:tbl>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading2"/>
  </w:pPr>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:t>Creazione nuova categoria</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading2"/>
  </w:pPr>
  <w:r>
    <w:t>Modifica categoria</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading2"/>
  </w:pPr>
  <w:r>
    <w:t>Eliminazione categoria</w:t>
  </w:r>
</w:p>

</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$rng.InsertXML($xmlFrag)
